# Updated cryptos list — refresh the Price (column D) and Volume(1h)
# (column E) figures for the coinranking.com snapshot on sheet1.
#
# Every value in columns D/E is stored as literal text in the workbook
# (e.g. "65.991.57", "  +0.31%  ") rather than as a number, so we force
# each target cell to Text format before writing the new string and then
# strip the format back off (Style = "Normal") so we don't leave behind
# a visible formatting change that isn't part of the data refresh.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = "D2";  Value = "66.014.62" },
    @{ Cell = "E2";  Value = "  +0.42%  " },
    @{ Cell = "D3";  Value = "3.315.63" },
    @{ Cell = "E3";  Value = "  +1.22%  " },
    @{ Cell = "E4";  Value = "  -0.06%  " },
    @{ Cell = "D5";  Value = "562.50" },
    @{ Cell = "E5";  Value = "  +1.12%  " },
    @{ Cell = "D6";  Value = "185.22" },
    @{ Cell = "E6";  Value = "  +1.04%  " },
    @{ Cell = "D7";  Value = "1.00" },
    @{ Cell = "E7";  Value = "  +0.12%  " },
    @{ Cell = "D8";  Value = "3.310.30" },
    @{ Cell = "E8";  Value = "  +1.33%  " },
    @{ Cell = "E9";  Value = "  -2.51%  " },
    @{ Cell = "E10"; Value = "  -5.18%  " },
    @{ Cell = "E11"; Value = "  -1.67%  " },
    @{ Cell = "D12"; Value = "45.87" },
    @{ Cell = "E12"; Value = "  -2.85%  " },
    @{ Cell = "E13"; Value = "  -1.08%  " },
    @{ Cell = "D14"; Value = "3.849.06" },
    @{ Cell = "E14"; Value = "  +1.42%  " },
    @{ Cell = "D15"; Value = "8.45" },
    @{ Cell = "E15"; Value = "  -2.33%  " },
    @{ Cell = "D16"; Value = "588.23" },
    @{ Cell = "E16"; Value = "  -9.12%  " },
    @{ Cell = "D17"; Value = "65.999.75" },
    @{ Cell = "E17"; Value = "  +0.37%  " },
    @{ Cell = "D19"; Value = "3.313.21" },
    @{ Cell = "E19"; Value = "  +1.15%  " },
    @{ Cell = "D20"; Value = "17.69" },
    @{ Cell = "E20"; Value = "  -1.82%  " },
    @{ Cell = "D21"; Value = "10.89" },
    @{ Cell = "E21"; Value = "  -3.91%  " },
    @{ Cell = "D23"; Value = "17.85" },
    @{ Cell = "E23"; Value = "  -2.42%  " },
    @{ Cell = "D24"; Value = "5.04" },
    @{ Cell = "E24"; Value = "  +2.27%  " },
    @{ Cell = "D25"; Value = "97.76" },
    @{ Cell = "E25"; Value = "  -10.00%  " },
    @{ Cell = "D26"; Value = "3.98" },
    @{ Cell = "E26"; Value = "  +0.21%  " },
    @{ Cell = "E27"; Value = "  +0.93%  " },
    @{ Cell = "E28"; Value = "  -2.58%  " },
    @{ Cell = "D29"; Value = "8.43" },
    @{ Cell = "E29"; Value = "  -2.52%  " },
    @{ Cell = "E30"; Value = "  +1.12%  " },
    @{ Cell = "D31"; Value = "6.63" },
    @{ Cell = "E31"; Value = "  +5.57%  " },
    @{ Cell = "D32"; Value = "565.55" },
    @{ Cell = "E32"; Value = "  +9.01%  " },
    @{ Cell = "D33"; Value = "3.67" },
    @{ Cell = "E33"; Value = "  -7.01%  " },
    @{ Cell = "D34"; Value = "10.81" },
    @{ Cell = "E34"; Value = "  -2.17%  " },
    @{ Cell = "D35"; Value = "3.773.41" },
    @{ Cell = "E35"; Value = "  +0.24%  " },
    @{ Cell = "E36"; Value = "  -1.71%  " },
    @{ Cell = "D37"; Value = "0.999" },
    @{ Cell = "E37"; Value = "  -0.12%  " },
    @{ Cell = "D38"; Value = "55.76" },
    @{ Cell = "E38"; Value = "  -2.96%  " },
    @{ Cell = "D39"; Value = "33.24" },
    @{ Cell = "E39"; Value = "  +1.12%  " },
    @{ Cell = "D40"; Value = "0.127" },
    @{ Cell = "E40"; Value = "  -2.81%  " },
    @{ Cell = "D41"; Value = "3.14" },
    @{ Cell = "E41"; Value = "  -7.77%  " },
    @{ Cell = "D42"; Value = "0.0₃0684" },
    @{ Cell = "E42"; Value = "  -6.84%  " },
    @{ Cell = "E43"; Value = "  +4.68%  " },
    @{ Cell = "E44"; Value = "  -5.69%  " },
    @{ Cell = "E45"; Value = "  -1.26%  " },
    @{ Cell = "D47"; Value = "3.06" },
    @{ Cell = "E47"; Value = "  -10.24%  " },
    @{ Cell = "E48"; Value = "  -2.31%  " },
    @{ Cell = "E49"; Value = "  +0.03%  " },
    @{ Cell = "E50"; Value = "  -3.32%  " },
    @{ Cell = "D51"; Value = "127.85" },
    @{ Cell = "E51"; Value = "  +5.00%  " }
)

foreach ($update in $updates) {
    $rng = $ws.Range($update.Cell)
    $rng.NumberFormat = "@"
    $rng.Value = $update.Value
    $rng.Style = "Normal"
}
